$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - copy formatting (style) from the existing H1 header cell
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy()
$ws.Range("J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J1").Value = "IF"

# Data rows for columns I and J
$data = @{
    2  = @(1, 4)
    3  = @(1, 5)
    4  = @(1, 4)
    5  = @(1, 2)
    6  = @(1, 3)
    7  = @(5, 7)
    8  = @(1, 4)
    9  = @(4, 5)
    10 = @(1, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
